$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Timesheet")

# --- Row 2 : 2026-01-19 ---
$ws.Range("B2").Value = "Keevil"
$ws.Range("C2").Value = 8
$ws.Range("E2").Value = 65
$ws.Range("F2").Value = 520

# --- Row 3 : 2026-01-20 ---
$ws.Range("B3").Value = "McGill"
$ws.Range("C3").Value = 9.5
$ws.Range("E3").Value = 65
$ws.Range("F3").Value = 617.5

# --- Row 4 : 2026-01-21 ---
$ws.Range("B4").Value = "Smith"
$ws.Range("C4").Value = 9
$ws.Range("E4").Value = 65
$ws.Range("F4").Value = 585

# --- Row 5 : 2026-01-22 ---
$ws.Range("B5").Value = "Bottomley"
$ws.Range("C5").Value = 9
$ws.Range("E5").Value = 65
$ws.Range("F5").Value = 585

# --- Row 6 : 2026-01-23 (Regular) ---
$ws.Range("B6").Value = "Varricchio"
$ws.Range("C6").Value = 4.5
$ws.Range("E6").Value = 65
$ws.Range("F6").Value = 292.5

# --- Row 7 : 2026-01-23 (OT) ---
$ws.Range("B7").Value = "Varricchio"
$ws.Range("C7").Value = 3.5
$ws.Range("E7").Value = 65
$ws.Range("F7").Value = 341.25

# --- Row 9 : SUBTOTAL ---
$ws.Range("C9").Value = 43.5
$ws.Range("D9").Value = "Reg: 40 / OT: 3.5"
$ws.Range("F9").Value = 2941.25
